$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "dhskah"
$ws.Range("H9").Value = "kjskn"
$ws.Range("L14").Value = "kjsnkjn"
[void]$ws.Range("L14").Select()
